$d = $word.ActiveDocument

# --- Location 1 -----------------------------------------------------------
# "...se cuenta con todos los componentes y acciones que debe de seguir
#  dichas aplicaciones del sistema" (unique phrasing, third occurrence of
#  "componentes" in the document) becomes "...objetos..." and gains a new
#  trailing clause " para llegar a su propósito". The (empty) "_GoBack"
#  bookmark now lives right after this new text, before the closing ".".
$old1 = "todos los componentes y acciones que debe de seguir dichas aplicaciones del sistema"
$new1 = "todos los objetos y acciones que debe de seguir dichas aplicaciones del sistema para llegar a su propósito"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# Re-locate the text we just wrote so a bookmark can be dropped right
# after it (immediately before the following "." run).
$rng1 = $d.Content
$rng1.Find.Execute($new1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $d.Range($rng1.End, $rng1.End)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Location 2 -------------------------------------------------------------
# The old "_GoBack" bookmark used to sit between "Correcciones " and
# "/ Recomendaciones" inside "Observaciones / Notas de Correcciones
# / Recomendaciones". Re-write that trailing part so the stale bookmark is
# dropped and the text reads as a single " / Recomendaciones" tail.
$old2 = "Correcciones / Recomendaciones"
$new2 = "Correcciones / Recomendaciones"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
